$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. New column P ("cruise missiles"): width + header style (new bold font)
# ------------------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 17.666666666666668   # renders as 18.5 in OOXML

$ws.Range("P1").Value = "cruise missiles"
$ws.Range("P1").Font.Bold = $true
$ws.Range("P1").Font.Name = "Arial"
$ws.Range("P1").Font.Size = 10
$ws.Range("P1").Font.Color = 0
$ws.Range("P1").HorizontalAlignment = -4108
$ws.Range("P1").VerticalAlignment = -4108

# ------------------------------------------------------------------
# 2. Fill column P with 0 for the already existing rows (2-67);
#    rows 68 and 69 (new) get 84, handled later together with the
#    rest of their row data.
# ------------------------------------------------------------------
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
}

# ------------------------------------------------------------------
# 3. Append the new daily rows (54-69).
#    Source data (columns B..P) for each new row, in order.
# ------------------------------------------------------------------
$newRows = @(
    @(20300,773,2002,376,127,66,165,146,1471,8,76,148,27,4,0),
    @(20600,790,2041,381,130,67,167,147,1487,8,76,155,27,4,0),
    @(20800,802,2063,386,132,67,169,150,1495,8,76,158,27,4,0),
    @(20900,815,2087,391,136,67,171,150,1504,8,76,165,27,4,0),
    @(21000,829,2118,393,136,67,172,151,1508,8,76,166,27,4,0),
    @(21200,838,2162,397,138,69,176,153,1523,8,76,172,27,4,0),
    @(21600,854,2205,403,143,69,177,154,1543,8,76,182,27,4,0),
    @(21800,873,2238,408,147,69,179,154,1557,8,76,191,28,4,0),
    @(21900,884,2258,411,149,69,181,154,1566,8,76,201,28,4,0),
    @(22100,918,2308,416,149,69,184,154,1643,8,76,205,31,4,0),
    @(22400,939,2342,421,149,71,185,155,1666,8,76,207,31,4,0),
    @(22800,970,2389,431,151,72,187,155,1688,8,76,215,31,4,0),
    @(23000,986,2418,435,151,73,189,155,1695,8,76,229,31,4,0),
    @(23200,1008,2445,436,151,77,190,155,1701,8,76,232,32,4,0),
    @(23500,1026,2471,451,151,80,192,155,1796,8,76,245,32,4,84),
    @(23800,1048,2519,459,152,80,194,155,1824,8,76,271,38,4,84)
)

$startRow = 54
$endRow = 69

# Clone number formats / fonts / alignment from the last existing row for
# the whole new block at once, so every new row matches the table's
# existing look (date style on A, plain style on B..O).
$ws.Range("A53:O53").Copy()
$ws.Range(("A" + $startRow + ":O" + $endRow)).PasteSpecial(-4122)
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Range("A" + $r).RowHeight = 15.75
}

# Running date formula for the whole new block in one shot, so Excel keeps
# it as a single shared formula (like the rest of column A).
$ws.Range(("A" + $startRow + ":A" + $endRow)).Formula = "=A53+1"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $vals[$c]
    }
}

# ------------------------------------------------------------------
# 4. Final selection / scroll position, mirroring the author's cursor.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 13
$ws.Range("P70").Select()
